$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6060816.5
$ws.Range("I33").Value = 221.66667
$ws.Range("J33").Value = 30303196
$ws.Range("K33").Value = 221.66667
$ws.Range("L33").Value = 30303196
$ws.Range("M33").Value = 7.333329999999989
$ws.Range("N33").Value = -30303654
$ws.Range("H38").Value = 149
$ws.Range("I38").Value = 149
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 447
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -75
$ws.Range("N38").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H106").Value = 78435960
$ws.Range("I106").Value = 23814024
$ws.Range("K106").Value = 23814024
$ws.Range("M106").Value = -23813393
$ws.Range("H132").Value = 2338.8667
$ws.Range("I132").Value = 2470.9285
$ws.Range("J132").Value = 490
$ws.Range("K132").Value = 7412.7855
$ws.Range("L132").Value = 1470
$ws.Range("M132").Value = -4882.7855
$ws.Range("N132").Value = -6530
$ws.Range("H135").Value = 2402.2
$ws.Range("I135").Value = 2507.4375
$ws.Range("J135").Value = 1981.25
$ws.Range("K135").Value = 22566.9375
$ws.Range("L135").Value = 17831.25
$ws.Range("M135").Value = -20031.9375
$ws.Range("N135").Value = -22901.25
$ws.Range("H137").Value = 1205.7539
$ws.Range("I137").Value = 781.6316
$ws.Range("J137").Value = 1802.6666
$ws.Range("K137").Value = 2344.8948
$ws.Range("L137").Value = 5407.9998
$ws.Range("M137").Value = 205.1052
$ws.Range("N137").Value = -10507.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 20300
$ws.Range("J24").Value = 20300
$ws.Range("L24").Value = 20300
$ws.Range("N24").Value = -21048
$ws.Range("H32").Value = 5715.4614
$ws.Range("I32").Value = 5461.729
$ws.Range("J32").Value = 6431.8823
$ws.Range("K32").Value = 5461.729
$ws.Range("L32").Value = 6431.8823
$ws.Range("M32").Value = -5174.729
$ws.Range("N32").Value = -7005.8823
$ws.Range("H61").Value = 4035.7673
$ws.Range("J61").Value = 1123.6471
$ws.Range("L61").Value = 1123.6471
$ws.Range("N61").Value = -1547.6471
$ws.Range("H100").Value = 20300
$ws.Range("J100").Value = 20300
$ws.Range("L100").Value = 20300
$ws.Range("N100").Value = -22464
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 3020.5454
$ws.Range("I132").Value = 1492.84
$ws.Range("J132").Value = 7794.625
$ws.Range("K132").Value = 4478.52
$ws.Range("L132").Value = 23383.875
$ws.Range("M132").Value = -1948.52
$ws.Range("N132").Value = -28443.875
$ws.Range("H134").Value = 19000
$ws.Range("J134").Value = 19000
$ws.Range("L134").Value = 19000
$ws.Range("N134").Value = -29140
$ws.Range("H136").Value = 4035.7673
$ws.Range("J136").Value = 1123.6471
$ws.Range("L136").Value = 3370.9413
$ws.Range("N136").Value = -8470.941299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3017619.5
$ws.Range("I7").Value = 4507626
$ws.Range("K7").Value = 4507626
$ws.Range("M7").Value = -4507513
$ws.Range("H138").Value = 59800
$ws.Range("J138").Value = 59800
$ws.Range("L138").Value = 59800
$ws.Range("N138").Value = -70080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8547759
$ws.Range("I16").Value = 10989717
$ws.Range("J16").Value = 906.5
$ws.Range("K16").Value = 10989717
$ws.Range("L16").Value = 906.5
$ws.Range("M16").Value = -10989430
$ws.Range("N16").Value = -1480.5
$ws.Range("H31").Value = 3777.3494
$ws.Range("I31").Value = 1639.44
$ws.Range("J31").Value = 4698.8623
$ws.Range("K31").Value = 1639.44
$ws.Range("L31").Value = 4698.8623
$ws.Range("M31").Value = -1344.44
$ws.Range("N31").Value = -5288.8623
$ws.Range("H34").Value = 3777.3494
$ws.Range("I34").Value = 1639.44
$ws.Range("J34").Value = 4698.8623
$ws.Range("K34").Value = 1639.44
$ws.Range("L34").Value = 4698.8623
$ws.Range("M34").Value = -1437.44
$ws.Range("N34").Value = -5102.8623
$ws.Range("H60").Value = 18750
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 18750
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 18750
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -19772
$ws.Range("H113").Value = 8547759
$ws.Range("I113").Value = 10989717
$ws.Range("J113").Value = 906.5
$ws.Range("K113").Value = 10989717
$ws.Range("L113").Value = 906.5
$ws.Range("M113").Value = -10987547
$ws.Range("N113").Value = -5246.5
$ws.Range("H132").Value = 2802
$ws.Range("I132").Value = 1785.1428
$ws.Range("J132").Value = 4225.6
$ws.Range("K132").Value = 5355.428400000001
$ws.Range("L132").Value = 12676.8
$ws.Range("M132").Value = -2825.428400000001
$ws.Range("N132").Value = -17736.8
$ws.Range("H135").Value = 33563.332
$ws.Range("J135").Value = 33563.332
$ws.Range("L135").Value = 33563.332
$ws.Range("N135").Value = -43703.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2556.6824
$ws.Range("I68").Value = 2939.4783
$ws.Range("J68").Value = 2105.1794
$ws.Range("K68").Value = 8818.4349
$ws.Range("L68").Value = 6315.5382
$ws.Range("M68").Value = -8007.4349
$ws.Range("N68").Value = -7937.5382
$ws.Range("H71").Value = 2556.6824
$ws.Range("I71").Value = 2939.4783
$ws.Range("J71").Value = 2105.1794
$ws.Range("K71").Value = 26455.3047
$ws.Range("L71").Value = 18946.6146
$ws.Range("M71").Value = -22399.3047
$ws.Range("N71").Value = -27058.6146
$ws.Range("H107").Value = 941.43665
$ws.Range("I107").Value = 309.3684
$ws.Range("J107").Value = 1172.3846
$ws.Range("K107").Value = 928.1052
$ws.Range("L107").Value = 3517.1538
$ws.Range("M107").Value = 991.8948
$ws.Range("N107").Value = -7357.1538
$ws.Range("H134").Value = 10430.942
$ws.Range("I134").Value = 11371.182
$ws.Range("K134").Value = 34113.546
$ws.Range("M134").Value = -29043.546
$ws.Range("H137").Value = 29426836
$ws.Range("I137").Value = 1657.4445
$ws.Range("J137").Value = 40019900
$ws.Range("K137").Value = 4972.333500000001
$ws.Range("L137").Value = 120059700
$ws.Range("M137").Value = 127.6664999999994
$ws.Range("N137").Value = -120069900
$ws.Range("H139").Value = 4288.3486
$ws.Range("I139").Value = 5840.4287
$ws.Range("J139").Value = 2806.818
$ws.Range("K139").Value = 17521.2861
$ws.Range("L139").Value = 8420.454000000002
$ws.Range("M139").Value = -12381.2861
$ws.Range("N139").Value = -18700.454
$ws.Range("H141").Value = 14624.19
$ws.Range("I141").Value = 10708.308
$ws.Range("J141").Value = 20987.5
$ws.Range("K141").Value = 32124.924
$ws.Range("L141").Value = 62962.5
$ws.Range("M141").Value = -26944.924
$ws.Range("N141").Value = -73322.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3242837
$ws.Range("I122").Value = 3242837
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9728511
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9726061
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 7326.222
$ws.Range("I126").Value = 8398.134
$ws.Range("J126").Value = 1966.6666
$ws.Range("K126").Value = 25194.402
$ws.Range("L126").Value = 5899.9998
$ws.Range("M126").Value = -22724.402
$ws.Range("N126").Value = -10839.9998
$ws.Range("H132").Value = 3485.4707
$ws.Range("I132").Value = 3545.6667
$ws.Range("J132").Value = 3452.6365
$ws.Range("K132").Value = 10637.0001
$ws.Range("L132").Value = 10357.9095
$ws.Range("M132").Value = -8107.000100000001
$ws.Range("N132").Value = -15417.9095
$ws.Range("H133").Value = 68545
$ws.Range("J133").Value = 68545
$ws.Range("L133").Value = 68545
$ws.Range("N133").Value = -78665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 25642056
$ws.Range("I46").Value = 30303920
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 30303920
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -30303732
$ws.Range("N46").Value = -2176
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H132").Value = 16674017
$ws.Range("I132").Value = 24083034
$ws.Range("J132").Value = 3726
$ws.Range("K132").Value = 72249102
$ws.Range("L132").Value = 11178
$ws.Range("M132").Value = -72246572
$ws.Range("N132").Value = -16238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1835.3
$ws.Range("I132").Value = 1692.4482
$ws.Range("J132").Value = 2211.9092
$ws.Range("K132").Value = 5077.3446
$ws.Range("L132").Value = 6635.7276
$ws.Range("M132").Value = -2547.3446
$ws.Range("N132").Value = -11695.7276
$ws.Range("H136").Value = 1070.8823
$ws.Range("I136").Value = 651.875
$ws.Range("J136").Value = 2076.5
$ws.Range("K136").Value = 2076.5
$ws.Range("L136").Value = 6229.5
$ws.Range("M136").Value = 594.375
$ws.Range("N136").Value = -11329.5
